# issue #5: stock data from json to db
# Adds three columns (category, source_file, index) to the "股票" (stock)
# worksheet, mirroring the extra metadata columns that the json->db
# pipeline now emits for every property record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Last row of data (header is row 1, data starts row 2).
$lastRow = 33

# Insert a new "category" column between "property_category" (H) and
# "date" (old I, now shifts to J). This naturally carries the existing
# "data row" cell style along with it, just like Excel's own
# Insert-Shift-Right behaviour.
$ws.Columns.Item(9).Insert()
$ws.Range("I1").Value = "category"
$ws.Range("I2:I" + $lastRow).Value = "normal"

# Append two more columns ("source_file" and "index") right after the
# existing last column (old K/legislator_id, now L). Inserting columns
# here (rather than just writing past the used range) keeps the same
# cell formatting the rest of the data rows use.
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(14).Insert()

$ws.Range("M1").Value = "source_file"
$ws.Range("M2:M" + $lastRow).Value = "tmpf37d1"

$ws.Range("N1").Value = "index"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value2
}
